# "Error Calculations and Plots"
# This workbook holds a data table (ID + columns A-F[named A,B,C,D,F in header])
# with some cells intentionally blanked out to simulate missing data.
# The edit removes two whole data rows ("RM 232" and "SC 92") - causing every
# row below them to shift up - and then re-randomizes which cells in the
# remaining rows are blanked out vs populated with a value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Step 1: remove the two rows that disappear from the table entirely.
# Row 26 is "RM 232"; after it is deleted, the old row 28 "SC 92" becomes
# row 27, so we delete that same (now-shifted) row index a second time.
# ---------------------------------------------------------------------------
$ws.Rows(26).Delete()
$ws.Rows(27).Delete()

# ---------------------------------------------------------------------------
# Step 2: after the rows above have shifted up, adjust the individual cells
# that differ between the (row-shifted) original data and the target data -
# some previously-missing cells now contain a value, and some previously
# populated cells are now blanked out.
# ---------------------------------------------------------------------------

# Row 2 (RM 2)
$ws.Cells.Item(2,5).ClearContents()

# Row 4 (RM 9)
$ws.Cells.Item(4,6).ClearContents()

# Row 5 (RM 14)
$ws.Cells.Item(5,5).Value = -5

# Row 6 (RM 21)
$ws.Cells.Item(6,4).Value = -14.2
$ws.Cells.Item(6,5).Value = -5.7
$ws.Cells.Item(6,6).Value = 16.43

# Row 8 (RM 38)
$ws.Cells.Item(8,4).ClearContents()

# Row 9 (RM 42)
$ws.Cells.Item(9,5).ClearContents()

# Row 10 (RM 52 a)
$ws.Cells.Item(10,5).ClearContents()

# Row 11 (RM 58)
$ws.Cells.Item(11,6).Value = 17.65

# Row 12 (RM 81)
$ws.Cells.Item(12,4).Value = -14.1
$ws.Cells.Item(12,6).ClearContents()

# Row 14 (RM 90)
$ws.Cells.Item(14,4).ClearContents()
$ws.Cells.Item(14,6).Value = 17.76

# Row 17 (RM 116)
$ws.Cells.Item(17,4).Value = -14.7
$ws.Cells.Item(17,6).ClearContents()

# Row 18 (RM 120)
$ws.Cells.Item(18,4).Value = -15.2

# Row 19 (RM 125)
$ws.Cells.Item(19,4).ClearContents()
$ws.Cells.Item(19,6).Value = 17.81

# Row 20 (RM 134)
$ws.Cells.Item(20,4).ClearContents()

# Row 21 (RM 135)
$ws.Cells.Item(21,6).Value = 16.58

# Row 22 (RM 138)
$ws.Cells.Item(22,6).Value = 16.81

# Row 23 (RM 140)
$ws.Cells.Item(23,4).Value = -13.9

# Row 24 (RM 142a)
$ws.Cells.Item(24,5).Value = -8.1

# Row 25 (RM 145)
$ws.Cells.Item(25,6).ClearContents()

# Row 26 (SC 5, was row 27 before the deletes)
$ws.Cells.Item(26,6).ClearContents()

# Row 27 (SC 101, was row 29 before the deletes)
$ws.Cells.Item(27,3).Value = 10
$ws.Cells.Item(27,4).ClearContents()
$ws.Cells.Item(27,6).ClearContents()

# Row 28 (SC 105, was row 30 before the deletes)
$ws.Cells.Item(28,3).ClearContents()
$ws.Cells.Item(28,5).ClearContents()
$ws.Cells.Item(28,6).ClearContents()

# Row 29 (SC 119, was row 31 before the deletes)
$ws.Cells.Item(29,3).ClearContents()

# Row 30 (SC 120, was row 32 before the deletes)
$ws.Cells.Item(30,3).Value = 11.4
$ws.Cells.Item(30,5).Value = -5.7

# Row 31 (SC 132, was row 33 before the deletes)
$ws.Cells.Item(31,6).Value = 17.18

# Row 32 (SC 193, was row 34 before the deletes)
$ws.Cells.Item(32,3).ClearContents()
